$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new mapping row: coal_coke -> coal_coke
$ws.Range("A6").Value = "coal_coke"
$ws.Range("B6").Value = "coal_coke"

# Update the selected cell to mimic the saved selection state (A7)
$ws.Range("A7").Select()
